# Auto update Excel log
# Appends new sensor-log rows to the "PIR" sheet (rows 121-133) and the
# "Humidity" sheet (rows 79-88), matching the source logger's plain-text
# cell format (dates/percentages stored as literal text, not auto-converted
# numbers/dates).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper data: PIR (motion) sheet new rows
# Columns: Date, Timestamp, Hour, Location, Value, Status
# ---------------------------------------------------------------------------
$pirRows = @(
    @{ Row=121; Timestamp="18:28:36" },
    @{ Row=122; Timestamp="18:28:39" },
    @{ Row=123; Timestamp="18:28:44" },
    @{ Row=124; Timestamp="18:28:49" },
    @{ Row=125; Timestamp="18:28:54" },
    @{ Row=126; Timestamp="18:28:59" },
    @{ Row=127; Timestamp="18:29:04" },
    @{ Row=128; Timestamp="18:29:09" },
    @{ Row=129; Timestamp="18:29:14" },
    @{ Row=130; Timestamp="18:29:19" },
    @{ Row=131; Timestamp="18:29:24" },
    @{ Row=132; Timestamp="18:29:29" },
    @{ Row=133; Timestamp="18:29:34" }
)

$ws = $wb.Worksheets.Item("PIR")

foreach ($entry in $pirRows) {
    $r = $entry.Row

    # Column A ("Date") looks like a date ("2026-01-30") and would be
    # auto-parsed into a date serial by normal cell entry, so force the
    # cell to Text first, matching how the rest of the log stores it.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = "2026-01-30"

    $ws.Cells.Item($r, 2).Value = $entry.Timestamp
    $ws.Cells.Item($r, 3).Value = "18:00"
    $ws.Cells.Item($r, 4).Value = "Bathroom"
    $ws.Cells.Item($r, 5).Value = "No Motion"
    $ws.Cells.Item($r, 6).Value = "Inactive"
}

# ---------------------------------------------------------------------------
# Helper data: Humidity sheet new rows
# Columns: Date, Timestamp, Hour, Location, Value, Status
# ---------------------------------------------------------------------------
$humidityRows = @(
    @{ Row=79; Timestamp="18:28:36"; Value="86.5%" },
    @{ Row=80; Timestamp="18:28:44"; Value="86.4%" },
    @{ Row=81; Timestamp="18:28:49"; Value="86.4%" },
    @{ Row=82; Timestamp="18:28:54"; Value="86.4%" },
    @{ Row=83; Timestamp="18:29:04"; Value="86.4%" },
    @{ Row=84; Timestamp="18:29:09"; Value="86.3%" },
    @{ Row=85; Timestamp="18:29:14"; Value="86.4%" },
    @{ Row=86; Timestamp="18:29:25"; Value="86.5%" },
    @{ Row=87; Timestamp="18:29:29"; Value="86.4%" },
    @{ Row=88; Timestamp="18:29:35"; Value="86.4%" }
)

$ws2 = $wb.Worksheets.Item("Humidity")

foreach ($entry in $humidityRows) {
    $r = $entry.Row

    $ws2.Cells.Item($r, 1).NumberFormat = "@"
    $ws2.Cells.Item($r, 1).Value = "2026-01-30"

    $ws2.Cells.Item($r, 2).Value = $entry.Timestamp
    $ws2.Cells.Item($r, 3).Value = "18:00"
    $ws2.Cells.Item($r, 4).Value = "Bathroom"

    # Column E ("Value") looks like a percentage ("86.5%") and would be
    # auto-parsed into a numeric percentage by normal cell entry, so force
    # the cell to Text first, matching how the rest of the log stores it.
    $ws2.Cells.Item($r, 5).NumberFormat = "@"
    $ws2.Cells.Item($r, 5).Value = $entry.Value

    $ws2.Cells.Item($r, 6).Value = "Active"
}
